# Auto - Update data with bot!
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "R 개발/분석을 위한 안전한 모듈 관리 - box"
$ws.Range("E3").Value = "https://lumiamitie.github.io/dev/r-module-using-box/"

$ws.Range("D12").Value = "[핸즈온 머신러닝 2판], [머신 러닝 교과서 3판] 사이킷런 1.0 업데이트 완료"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/10/24/%ed%95%b8%ec%a6%88%ec%98%a8-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d-2%ed%8c%90-%eb%a8%b8%ec%8b%a0-%eb%9f%ac%eb%8b%9d-%ea%b5%90%ea%b3%bc%ec%84%9c-3%ed%8c%90-%ec%82%ac%ec%9d%b4%ed%82%b7%eb%9f%b0-1-0/"

$ws.Range("D39").Value = "Visualize your data with Facets"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Visualize-your-data-with-Facets-1"

$ws.Range("D41").Value = "AI에 대한 Tesla의 포부"
$ws.Range("E41").Value = "http://cloudinsight.net/ai/ai%ec%97%90-%eb%8c%80%ed%95%9c-tesla%ec%9d%98-%ed%8f%ac%eb%b6%80/"

$ws.Range("D44").Value = "Qualcomm Ventures Portfolio (2) - Augury"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/103"
